# UniformF-HW20.xlsx : add the "Holden" sampling scheme to the simulation
# comparison table.
#
# The sheet previously had a duplicated block of columns (B..T repeated as
# U..AD) and 18 scheme rows (rows 3-19). This edit:
#   1. Removes the duplicate column block (U:AD), which also normalizes the
#      sheet dimension / row spans.
#   2. Reorders the HKL-plane header labels in row 2 (C2:T2).
#   3. Renames the 4 "HexGrid-90degTilt*" rows (16-19) to the new
#      "Holden2.5/5/10/15" scheme names.
#   4. Appends 4 new rows (20-23) re-adding the "HexGrid-90degTilt*" schemes
#      at the bottom of the table, with the same C:T = 1 fill as every other
#      scheme row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the duplicate right-hand block of columns (U:AD). This also
#    updates <dimension> and each row's spans="" automatically.
$ws.Range("U1:AD19").EntireColumn.Delete()

# 2. Row 2 header labels (HKL planes), reordered.
$row2 = @("[3, 2, 1]", "[1, 1, 0]", "[2, 2, 2]", "[3, 1, 0]", "[2, 2, 0]", "[2, 0, 0]", "[2, 1, 1]", "[4, 0, 0]", "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique")
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])2").Value = $row2[$i]
}

# 3. Rename rows 16-19 from the HexGrid variants to the new Holden scheme.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 4. Append the displaced HexGrid rows at the bottom (rows 20-23), matching
#    the formatting (A column bold/bordered style) and the C:T = 1 fill used
#    by every other scheme row.
$newRows = @(
    @(20, 18, "HexGrid-90degTilt2.5degRes"),
    @(21, 19, "HexGrid-90degTilt5degRes"),
    @(22, 20, "HexGrid-90degTilt10degRes"),
    @(23, 21, "HexGrid-90degTilt15degRes")
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $aVal = $entry[1]
    $bVal = $entry[2]

    $ws.Range("A19").Copy($ws.Range("A$r"))
    $ws.Range("A$r").Value = $aVal
    $ws.Range("B$r").Value = $bVal

    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = 1
    }
}
